# Fruta / hortaliza, semanal
# Inserts one new weekly price record into the "Naranja" sheet at row 189,
# pushing the existing rows 189:291 down to 190:292 (dimension becomes
# A1:T292).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 189:291 down one row, creating a blank row 189 for the new record.
$ws.Rows(189).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(189, 1).Value  = 4
$ws.Cells.Item(189, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(189, 3).Value  = "Los Lagos"
$ws.Cells.Item(189, 4).Value  = 44518
$ws.Cells.Item(189, 5).Value  = 10
$ws.Cells.Item(189, 6).Value  = "Fruta"
$ws.Cells.Item(189, 7).Value  = 100102
$ws.Cells.Item(189, 8).Value  = "Cítricos"
$ws.Cells.Item(189, 9).Value  = 100102005
$ws.Cells.Item(189, 10).Value = "Naranja"
$ws.Cells.Item(189, 11).Value = "Navel Late"
$ws.Cells.Item(189, 12).Value = "Primera"
$ws.Cells.Item(189, 13).Value = 400
$ws.Cells.Item(189, 14).Value = 12000
$ws.Cells.Item(189, 15).Value = 13000
$ws.Cells.Item(189, 16).Value = 12500
$ws.Cells.Item(189, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(189, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(189, 19).Value = 833
$ws.Cells.Item(189, 20).Value = 15
